$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 36,2
$data[0,0] = "https://www.udemy.com/course/full-stack-programming-for-complete-beginners-in-python/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-gIbnrEzx9KRq46QVSLoN5w&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEAGAIN2"
$data[0,1] = $false
$data[1,0] = "https://www.udemy.com/course/bootstrap-5-with-5-projects-in-hindi-urdu/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-FvQ_.CTp0rL1vTOF2Xk75w&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEBOOTSTRAP5"
$data[1,1] = $false
$data[2,0] = "https://www.udemy.com/course/jenkins-github-and-aws-in-practice/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-T9l1LjEiqk1DtyQUZ0MRDQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=TAPIR4U"
$data[2,1] = $false
$data[3,0] = "https://www.udemy.com/course/nlp-with-transformers/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-8_TKEVmD148SLcYn9Lwr0w&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREE43"
$data[3,1] = $false
$data[4,0] = "https://www.udemy.com/course/fundamental-question-on-industrial-electronics/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-2jcaES8HGbMyKcOoyh4owA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=03476819394E431897AD"
$data[4,1] = $false
$data[5,0] = "https://www.udemy.com/course/certified-kubernetes-administrator-cka-practice-exams-2021-g/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-KftDcNRdbkLXV7Y.W89ouQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEDAYS"
$data[5,1] = $false
$data[6,0] = "https://www.udemy.com/course/aws-certified-security-specialty-practice-exams-latest-2021/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-mS4cQT2kgrkwCTWPhJf3pw&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEDAYS"
$data[6,1] = $false
$data[7,0] = "https://www.udemy.com/course/introduction-to-forex-learn-to-trade-forex-by-yourself/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-YlmLkeHLGgDClMVe_x35Mg&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=90B03C13055C47A28233"
$data[7,1] = $false
$data[8,0] = "https://www.udemy.com/course/the-complete-introduction-to-the-deep-web/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-ejoNkDbidMtJkVt3Dik03g&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=BFEF8CEA0D5C02B295EC"
$data[8,1] = $false
$data[9,0] = "https://www.udemy.com/course/java-programming-complete-beginner-to-advanced/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-kFna2ZPAJFIvysrg1IHSXQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=657F843318F537"
$data[9,1] = $false
$data[10,0] = "https://www.udemy.com/course/learn-guitar/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-7v0f5KRpF9pOzWBLcs4awg&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEMAY2021"
$data[10,1] = $false
$data[11,0] = "https://www.udemy.com/course/mathematics-software-development/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-4bG5mVFPyenfGx9DIxxGqg&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=716CB1EDFA7A3BDBA012"
$data[11,1] = $false
$data[12,0] = "https://www.udemy.com/course/learn-html5-in-depth-with-real-world-examples/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-TiBjtkhZWAdnGjFKUfXeYA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREE-JUNE"
$data[12,1] = $false
$data[13,0] = "https://www.udemy.com/course/five-proven-steps-to-real-estate-investing-success/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-xiCN6r1IgPMk4R7tTpCC0A&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FPSREEXPJUNE72021"
$data[13,1] = $false
$data[14,0] = "https://www.udemy.com/course/css3-in-hindi/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-8glOs2QMsLJlDbNfDfnc1A&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEHTML"
$data[14,1] = $false
$data[15,0] = "https://www.udemy.com/course/the-complete-nft-non-fungible-tokens-course-for-artists/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-z8Qt6L.JD4drRJYzQBiIbA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=C31651B940736748184A"
$data[15,1] = $false
$data[16,0] = "https://www.udemy.com/course/criminology-fundamentals-of-criminal-psychology-and-law/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-_w5DXJ__glaVUmoWmPyr_g&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=TECHBINZ"
$data[16,1] = $false
$data[17,0] = "https://www.udemy.com/course/malware-analysis-of-documents/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-M60CuswLMGrVodxClHsuyw&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FMADFREE_JUN2"
$data[17,1] = $false
$data[18,0] = "https://www.udemy.com/course/togaf-9-practice-exams-2021-combined-level-1-and-2/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-V6Bt8vkXdVJ_Rt7f.0ZX0A&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEDAYS"
$data[18,1] = $false
$data[19,0] = "https://www.udemy.com/course/nestjs-zero-to-hero/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-vuix.OYJdweFzxNEc4xWJA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=JUNE_REPRODUCTION"
$data[19,1] = $false
$data[20,0] = "https://www.udemy.com/course/powerbi-hero/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-fSCXMHSvvdOYIAtMa0V.CQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEDASHBI"
$data[20,1] = $false
$data[21,0] = "https://www.udemy.com/course/pmp-practice-test-project-management-professional-2021/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-3NzgBtrhXZl306P5kOdkQQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=88B2BE4D7ADDB73DC2A8"
$data[21,1] = $false
$data[22,0] = "https://www.udemy.com/course/procreate-sketch-draw-and-paint-a-shirt-design/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-Ka.lQqaDDH5LeKmXA62xAA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=479DD8FB85D20883AD98"
$data[22,1] = $false
$data[23,0] = "https://www.udemy.com/course/certified-ethical-hacker-ceh-v11-practice-exams-new-2021-p/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-C_tWhI4fou94y_Fuhc.GTA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEDAYS"
$data[23,1] = $false
$data[24,0] = "https://www.udemy.com/course/statistics-with-r-beginner-level/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-wHtlPlYCyDJnSo1Ig1G7sQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=BLIZZARD"
$data[24,1] = $false
$data[25,0] = "https://www.udemy.com/course/introduction-to-quantum-computing/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-NwgU.ZGMVqjWDGDnRgRf8w&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=1E58CF6234A770ADB7D4"
$data[25,1] = $false
$data[26,0] = "https://www.udemy.com/course/internet-and-web-development-fundamentals/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-JA4ZkVqndBxpGFa4FQ_yhA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=YOUACCEL44184"
$data[26,1] = $false
$data[27,0] = "https://www.udemy.com/course/learn-to-create-ai-voice-assistant-jarvis-with-python/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-XFOMo1r_V.kB_SehUO8kvw&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=5B0D13596E86F8590815"
$data[27,1] = $false
$data[28,0] = "https://www.udemy.com/course/advanced-neural-networks-in-r-a-practical-approach/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-D1Y2FQMwQso0jlxOpoolDQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=BLIZZARD"
$data[28,1] = $false
$data[29,0] = "https://www.udemy.com/course/master-complete-statistics-for-computer-science-i/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-5zpyg0CV2g64pUfJxfd3RQ&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEMCSCS8"
$data[29,1] = $false
$data[30,0] = "https://www.udemy.com/course/automate/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-2RAEBfwHARDAqISmMGjr6g&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=JUN2021FREE"
$data[30,1] = $false
$data[31,0] = "https://www.udemy.com/course/profitable-binary-trading-beginners-guide/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-zJCCLL3oT5a2984mJvElkw&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=F38EAFE4B41DC5727483"
$data[31,1] = $false
$data[32,0] = "https://www.udemy.com/course/complete-progressive-web-app-bootcamp/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-A9qitzgGelGw6tDMleaNfw&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FD68EF3DA859515E0BA4"
$data[32,1] = $false
$data[33,0] = "https://www.udemy.com/course/best-sap-fico-tutorial-for-beginners-freshers/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-zM9AA4RsBoUM1evDQAly9Q&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=FREEJUNE"
$data[33,1] = $false
$data[34,0] = "https://www.udemy.com/course/applied-ethical-hacking-and-rules-of-engagement/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-uOqvxC2Uedsk38t.dec0nA&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=AFFD464DDCD4D7C2AA17"
$data[34,1] = $false
$data[35,0] = "https://www.udemy.com/course/complete-linux-command-line-and-terminal-productivity/?ranMID=39197&ranEAID=%2F7fFXpljNdk&ranSiteID=_7fFXpljNdk-qNlUgfwBDe0d1W_LaYdP8A&LSNPUBID=%2F7fFXpljNdk&utm_source=aff-campaign&utm_medium=udemyads&couponCode=5CAEA2E49267B0A52E44"
$data[35,1] = $false

$ws.Range("A2:B37").Value2 = $data
